$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 5.124099999999999
$ws.Range("C2").Value = 21.80395
$ws.Range("F2").Value = 2035.1345
$ws.Range("G2").Value = 1832.0381
$ws.Range("H2").Value = 203.0962
$ws.Range("I2").Value = 203.0962
$ws.Range("K2").Value = 2180.4748
$ws.Range("L2").Value = 1832.041
$ws.Range("M2").Value = 348.4338
$ws.Range("N2").Value = 348.4338
$ws.Range("B3").Value = 5.138
$ws.Range("C3").Value = 21.881
$ws.Range("F3").Value = 2060.294
$ws.Range("G3").Value = 1834.905
$ws.Range("H3").Value = 225.389
$ws.Range("I3").Value = 225.389
$ws.Range("K3").Value = 2179.5908
$ws.Range("L3").Value = 1834.945
$ws.Range("M3").Value = 344.6458
$ws.Range("N3").Value = 344.6458
$ws.Range("B4").Value = 5.911
$ws.Range("C4").Value = 29.315
$ws.Range("F4").Value = 2896.988
$ws.Range("G4").Value = 2688.579
$ws.Range("H4").Value = 208.409
$ws.Range("I4").Value = 208.409
$ws.Range("K4").Value = 2212.8646
$ws.Range("L4").Value = 2049.397
$ws.Range("M4").Value = 163.4676
$ws.Range("N4").Value = 163.4676

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 5.151
$ws.Range("C2").Value = 22.1146
$ws.Range("F2").Value = 2054.9738
$ws.Range("G2").Value = 1840.3716
$ws.Range("H2").Value = 214.60205
$ws.Range("I2").Value = 214.60205
$ws.Range("K2").Value = 2088.428
$ws.Range("L2").Value = 1840.376
$ws.Range("M2").Value = 248.052
$ws.Range("N2").Value = 248.052
$ws.Range("B3").Value = 5.259
$ws.Range("C3").Value = 21.771
$ws.Range("F3").Value = 2082.455
$ws.Range("G3").Value = 1847.149
$ws.Range("H3").Value = 235.305
$ws.Range("I3").Value = 235.305
$ws.Range("K3").Value = 2089.224
$ws.Range("L3").Value = 1847.111
$ws.Range("M3").Value = 242.113
$ws.Range("N3").Value = 242.113
$ws.Range("B4").Value = 5.449
$ws.Range("C4").Value = 26.348
$ws.Range("F4").Value = 2897.878
$ws.Range("G4").Value = 2520.519
$ws.Range("H4").Value = 377.359
$ws.Range("I4").Value = 377.359
$ws.Range("K4").Value = 2102.4638
$ws.Range("L4").Value = 1946.03
$ws.Range("M4").Value = 156.4338
$ws.Range("N4").Value = 156.4338

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 5.0567
$ws.Range("C2").Value = 21.8227
$ws.Range("F2").Value = 2032.659
$ws.Range("G2").Value = 1824.5455
$ws.Range("H2").Value = 208.1134
$ws.Range("I2").Value = 208.1134
$ws.Range("K2").Value = 2197.0652
$ws.Range("L2").Value = 1824.549
$ws.Range("M2").Value = 372.5162
$ws.Range("N2").Value = 372.5162
$ws.Range("B3").Value = 5.145
$ws.Range("C3").Value = 21.772
$ws.Range("F3").Value = 2056.831
$ws.Range("G3").Value = 1833.964
$ws.Range("H3").Value = 222.867
$ws.Range("I3").Value = 222.867
$ws.Range("K3").Value = 2193.6706
$ws.Range("L3").Value = 1833.924
$ws.Range("M3").Value = 359.7466
$ws.Range("N3").Value = 359.7466
$ws.Range("B4").Value = 5.449
$ws.Range("C4").Value = 26.348
$ws.Range("F4").Value = 2897.878
$ws.Range("G4").Value = 2520.519
$ws.Range("H4").Value = 377.359
$ws.Range("I4").Value = 377.359
$ws.Range("K4").Value = 2186.239
$ws.Range("L4").Value = 1946.03
$ws.Range("M4").Value = 240.209
$ws.Range("N4").Value = 240.209

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 5.2129
$ws.Range("C2").Value = 22.9933
$ws.Range("F2").Value = 2099.346500000001
$ws.Range("G2").Value = 1862.30705
$ws.Range("H2").Value = 237.0393
$ws.Range("I2").Value = 237.0393
$ws.Range("K2").Value = 1945.4626
$ws.Range("L2").Value = 1862.311
$ws.Range("M2").Value = 83.1516
$ws.Range("N2").Value = 83.1516
$ws.Range("B3").Value = 5.336
$ws.Range("C3").Value = 22.678
$ws.Range("F3").Value = 2118.127
$ws.Range("G3").Value = 1871.255
$ws.Range("H3").Value = 246.872
$ws.Range("I3").Value = 246.872
$ws.Range("K3").Value = 1949.3042
$ws.Range("L3").Value = 1871.271
$ws.Range("M3").Value = 78.03319999999999
$ws.Range("N3").Value = 78.03319999999999
$ws.Range("B4").Value = 5.449
$ws.Range("C4").Value = 26.348
$ws.Range("F4").Value = 2897.878
$ws.Range("G4").Value = 2520.519
$ws.Range("H4").Value = 377.359
$ws.Range("I4").Value = 377.359
$ws.Range("K4").Value = 1991.9494
$ws.Range("L4").Value = 1946.03
$ws.Range("M4").Value = 45.9194
$ws.Range("N4").Value = 45.9194

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 5.1563
$ws.Range("C2").Value = 22.60825
$ws.Range("F2").Value = 2074.068400000001
$ws.Range("G2").Value = 1849.28095
$ws.Range("H2").Value = 224.78725
$ws.Range("I2").Value = 224.78725
$ws.Range("K2").Value = 2024.9946
$ws.Range("L2").Value = 1849.285
$ws.Range("M2").Value = 175.7096
$ws.Range("N2").Value = 175.7096
$ws.Range("B3").Value = 5.233
$ws.Range("C3").Value = 22.501
$ws.Range("F3").Value = 2098.47
$ws.Range("G3").Value = 1856.43
$ws.Range("H3").Value = 242.04
$ws.Range("I3").Value = 242.04
$ws.Range("K3").Value = 2024.8878
$ws.Range("L3").Value = 1856.367
$ws.Range("M3").Value = 168.5208
$ws.Range("N3").Value = 168.5208
$ws.Range("B4").Value = 5.449
$ws.Range("C4").Value = 26.348
$ws.Range("F4").Value = 2897.878
$ws.Range("G4").Value = 2520.519
$ws.Range("H4").Value = 377.359
$ws.Range("I4").Value = 377.359
$ws.Range("K4").Value = 2053.421
$ws.Range("L4").Value = 1946.03
$ws.Range("M4").Value = 107.391
$ws.Range("N4").Value = 107.391

Write-Output "Update complete"